# Bug fix: "finding the name of experience" — the CV-keyword extractor was
# mis-aligned by one slot (it was reading the *count* suffix instead of the full
# "<keyword> : <count>" phrase), which shifted every downstream keyword/score cell.
# This script re-keys Sheet1 ("nihad-azimli-resume__jodel_job") so the CV KEYWORDS
# (col E), CLUSTER MUST HAVE MATCH (col F) and the score cells (D2/G2/I2/K2) match
# the corrected extraction output, and appends the 6 extra keyword rows it now finds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: summary cells ---
# NOTE: D2, G2, I2 and K2 hold numeric-looking text ("24.69", "25.0", "0.0") that must
# be stored as TEXT (shared string), not as an Excel number. A leading single-quote
# is the standard Excel way to force text entry for a numeric-looking value.
$ws.Range("D2").Value = "'24.69"
$ws.Range("E2").Value = "data engineer : 1"
$ws.Range("F2").Value = "spark : 2"
$ws.Range("G2").Value = "'25.0"
$ws.Range("I2").Value = "'0.0"
$ws.Range("K2").Value = "'0.0"

# --- Rows 3-7: E (CV keyword) and F (matched cluster keyword) columns ---
$ws.Range("E3").Value = "engineer : 1"
$ws.Range("F3").Value = "python : 1"
$ws.Range("E4").Value = "apache : 4"
$ws.Range("F4").Value = "redshift : 3"
$ws.Range("E5").Value = "celery : 1"
$ws.Range("F5").Value = "engineering : 1"
$ws.Range("E6").Value = "python : 8"
$ws.Range("F6").Value = "s3 : 2"
$ws.Range("E7").Value = "aws : 16"
$ws.Range("F7").Value = "aws : 2"

# --- Rows 8-61: E column (rest of the original keyword list, re-aligned) ---
$ws.Range("E8").Value = "s3 : 2"
$ws.Range("E9").Value = "lambda : 3"
$ws.Range("E10").Value = "research : 1"
$ws.Range("E11").Value = "data migration : 1"
$ws.Range("E12").Value = "migration : 1"
$ws.Range("E13").Value = "mongodb : 5"
$ws.Range("E14").Value = "apache spark : 4"
$ws.Range("E15").Value = "spark : 2"
$ws.Range("E16").Value = "kubernetes : 4"
$ws.Range("E17").Value = "development : 1"
$ws.Range("E18").Value = "kinesis : 2"
$ws.Range("E19").Value = "dynamodb : 4"
$ws.Range("E20").Value = "reports : 1"
$ws.Range("E21").Value = "qlik : 4"
$ws.Range("E22").Value = "conversion : 1"
$ws.Range("E23").Value = "sql : 1"
$ws.Range("E24").Value = "ssis : 1"
$ws.Range("E25").Value = "etl : 1"
$ws.Range("E26").Value = "selenium : 2"
$ws.Range("E27").Value = "sql queries : 1"
$ws.Range("E28").Value = "queries : 1"
$ws.Range("E29").Value = "stored procedures : 1"
$ws.Range("E30").Value = "bi : 1"
$ws.Range("E31").Value = "developer : 3"
$ws.Range("E32").Value = "bamboo : 1"
$ws.Range("E33").Value = "deployment : 1"
$ws.Range("E34").Value = "software developer : 3"
$ws.Range("E35").Value = "software : 2"
$ws.Range("E36").Value = "analysis : 3"
$ws.Range("E37").Value = "project : 3"
$ws.Range("E38").Value = "data mining : 3"
$ws.Range("E39").Value = "mining : 2"
$ws.Range("E40").Value = "tensorflow : 3"
$ws.Range("E41").Value = "intern : 1"
$ws.Range("E42").Value = "computer engineering : 1"
$ws.Range("E43").Value = "engineering : 1"
$ws.Range("E44").Value = "electrical : 1"
$ws.Range("E45").Value = "electronics : 1"
$ws.Range("E46").Value = "ieee : 1"
$ws.Range("E47").Value = "hybrid : 1"
$ws.Range("E48").Value = "design : 1"
$ws.Range("E49").Value = "github : 2"
$ws.Range("E50").Value = "tools : 2"
$ws.Range("E51").Value = "docker : 2"
$ws.Range("E52").Value = "ansible : 2"
$ws.Range("E53").Value = "cloudformation : 1"
$ws.Range("E54").Value = "azure : 2"
$ws.Range("E55").Value = "amazon : 1"
$ws.Range("E56").Value = "web services : 1"
$ws.Range("E57").Value = "microsoft azure : 1"
$ws.Range("E58").Value = "databases : 1"
$ws.Range("E59").Value = "redshift : 2"
$ws.Range("E60").Value = "mysql : 2"
$ws.Range("E61").Value = "unix shell : 2"

# --- Rows 62-67: brand-new keyword rows found by the fixed extractor ---
$ws.Range("E62").Value = "shell : 1"
$ws.Range("E63").Value = "java : 1"
$ws.Range("E64").Value = "bi tools : 1"
$ws.Range("E65").Value = "tableau : 1"
$ws.Range("E66").Value = "english : 1"
$ws.Range("E67").Value = "amazon web services : 1"
